# The workbook has a single worksheet that needs to be renamed from
# "Collection_EG" to "CRF_EG" (the _xlnm._FilterDatabase defined name,
# which refers to this sheet, updates automatically when we rename it).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_EG")
$ws.Activate()
$ws.Name = "CRF_EG"

# Move the active selection on that sheet from R4 to L6.
$ws.Range("L6").Select()
